$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.518.54"
$ws.Range("E2").Value = "  +1.27%  "

$ws.Range("D3").Value = "3.386.82"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.52"
$ws.Range("E5").Value = "  +1.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.66"
$ws.Range("E6").Value = "  +2.15%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.385.62"
$ws.Range("E8").Value = "  +1.06%  "

$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("E11").Value = "  +3.19%  "

$ws.Range("E12").Value = "  +0.75%  "

$ws.Range("D13").Value = "3.961.36"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("E14").Value = "  +2.59%  "

$ws.Range("E15").Value = "  +3.06%  "

$ws.Range("D16").Value = "3.384.81"
$ws.Range("E16").Value = "  +1.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.78"
$ws.Range("E17").Value = "  +3.36%  "

$ws.Range("D18").Value = "61.563.51"
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("E19").Value = "  +2.37%  "

$ws.Range("E20").Value = "  +2.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.40"
$ws.Range("E21").Value = "  +0.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.97"
$ws.Range("E22").Value = "  +1.32%  "

$ws.Range("E23").Value = "  -2.26%  "

$ws.Range("D24").Value = "3.528.64"
$ws.Range("E24").Value = "  +1.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").Value = "  +8.79%  "

$ws.Range("E27").Value = "  +1.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.72"
$ws.Range("E28").Value = "  +4.54%  "

$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.162"
$ws.Range("E31").Value = "  +5.64%  "

$ws.Range("E32").Value = "  +2.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.17"
$ws.Range("E33").Value = "  +2.03%  "

$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.45"
$ws.Range("E35").Value = "  +0.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.30"
$ws.Range("E36").Value = "  -3.45%  "

$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.83"
$ws.Range("E38").Value = "  -0.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.64"
$ws.Range("E39").Value = "  +1.97%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("B41").Value = "ONDO"
$ws.Range("C41").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.24"
$ws.Range("E41").Value = "  +3.56%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.779"
$ws.Range("E43").Value = "  +3.35%  "

$ws.Range("E44").Value = "  +8.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.42"
$ws.Range("E45").Value = "  +1.10%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.46"
$ws.Range("E46").Value = "  +0.68%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.81"
$ws.Range("E47").Value = "  +10.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.84"
$ws.Range("E48").Value = "  -1.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.62"
$ws.Range("E49").Value = "  -2.14%  "

$ws.Range("D50").Value = "2.346.78"
$ws.Range("E50").Value = "  +5.82%  "

$ws.Range("E51").Value = "  -0.94%  "
